$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.018.03'
$ws.Range("E2").Value = '  +0.69%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.067.53'
$ws.Range("E3").Value = '  +0.17%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '537.15'
$ws.Range("E5").Value = '  -0.43%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.84'
$ws.Range("E6").Value = '  +2.72%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.059.38'
$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("E9").Value = '  +0.69%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.155'
$ws.Range("E10").Value = '  +1.02%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.19'
$ws.Range("E11").Value = '  +0.72%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.452'
$ws.Range("E12").Value = '  -2.01%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000221'
$ws.Range("E13").Value = '  +0.84%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.32'
$ws.Range("E14").Value = '  -1.04%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.563.87'
$ws.Range("E15").Value = '  +1.05%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.031.29'
$ws.Range("E16").Value = '  +0.66%  '

$ws.Range("E17").Value = '  +1.60%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.066.60'
$ws.Range("E18").Value = '  +0.02%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.61'
$ws.Range("E19").Value = '  -0.49%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '468.74'
$ws.Range("E20").Value = '  -2.14%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.39'
$ws.Range("E21").Value = '  -0.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.693'
$ws.Range("E22").Value = '  -2.06%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.98'
$ws.Range("E23").Value = '  -2.87%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.32'
$ws.Range("E24").Value = '  -0.34%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.08'
$ws.Range("E25").Value = '  +0.62%  '

$ws.Range("E26").Value = '  +0.09%  '

$ws.Range("E27").Value = '  -0.91%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.85'
$ws.Range("E28").Value = '  -4.46%  '

$ws.Range("E29").Value = '  -0.03%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '26.05'
$ws.Range("E30").Value = '  +0.07%  '

$ws.Range("E31").Value = '  +5.14%  '

$ws.Range("E32").Value = '  -2.89%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '58.94'
$ws.Range("E33").Value = '  +1.66%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.29'
$ws.Range("E34").Value = '  -5.08%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.45'
$ws.Range("E35").Value = '  +5.77%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.92'
$ws.Range("E36").Value = '  -0.72%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '479.90'
$ws.Range("E37").Value = '  -2.29%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.258.99'
$ws.Range("E38").Value = '  +3.96%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0395'
$ws.Range("E39").Value = '  +0.90%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0790'
$ws.Range("E40").Value = '  -0.68%  '

$ws.Range("E41").Value = '  +1.03%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.10'
$ws.Range("E42").Value = '  +0.63%  '

$ws.Range("E43").Value = '  +1.44%  '

$ws.Range("E44").Value = '  -1.27%  '

$ws.Range("E45").Value = '  +0.15%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '25.24'
$ws.Range("E46").Value = '  +2.56%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '122.98'
$ws.Range("E47").Value = '  +3.75%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.99'
$ws.Range("E48").Value = '  -2.38%  '

$ws.Range("E49").Value = '  +0.87%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₃0518'
$ws.Range("E50").Value = '  +2.31%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.00'
$ws.Range("E51").Value = '  -0.41%  '
